$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping-tipo-vivienda")

# Map of current (slug) values in column A to the new human-readable labels
$values = @(
    "Residencias de trabajadores",
    "Instituciones religiosas",
    "Vacías",
    "Internados, academias y escuelas militares, seminarios,¿",
    "Secundarias",
    "Convencionales",
    "Instituciones penitenciarias",
    "Asilos o residencias de ancianos",
    "Colegios mayores, residencias de estudiantes",
    "Otro tipo de colectivo",
    "Otro tipo",
    "Instituciones para personas con discapacidades",
    "Albergues para marginados sociales",
    "Hospitales de larga estancia",
    "Establecimientos militares",
    "Hoteles, pensiones, albergues,¿",
    "Otras instituciones de asistencias social a la infancia, juventud,¿",
    "Hospitales psiquiátricos",
    "Alojamientos",
    "Hospitales generales y especiales de corta estancia"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
